$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.405.38'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.834.64'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.69'
$ws.Range("E5").Value = '  -1.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5254'
$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2830'
$ws.Range("E8").Value = '  -11.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06906'
$ws.Range("E9").Value = '  +1.92%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.845.57'
$ws.Range("E10").Value = '  +0.45%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.21'
$ws.Range("E11").Value = '  -13.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.6958'
$ws.Range("E12").Value = '  -10.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07091'
$ws.Range("E13").Value = '  -8.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.35'
$ws.Range("E14").Value = '  -0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.853'
$ws.Range("E15").Value = '  -2.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.14'
$ws.Range("E18").Value = '  -4.84%  '

$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '26.440.76'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007187'
$ws.Range("E20").Value = '  -9.21%  '

$ws.Range("D21").Value = '2.086.69'
$ws.Range("E21").Value = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.487'
$ws.Range("E22").Value = '  -2.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.837'
$ws.Range("E23").Value = '  -1.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.894'
$ws.Range("E24").Value = '  -4.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.02'
$ws.Range("E25").Value = '  +0.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.681'
$ws.Range("E26").Value = '  -0.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.043'
$ws.Range("E27").Value = '  -6.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.57'
$ws.Range("E28").Value = '  -1.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '109.22'
$ws.Range("E29").Value = '  -2.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.082'
$ws.Range("E30").Value = '  -1.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08734'
$ws.Range("E31").Value = '  +0.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.832'
$ws.Range("E32").Value = '  -5.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04683'
$ws.Range("E33").Value = '  -3.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.908'
$ws.Range("E34").Value = '  +1.80%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.107'
$ws.Range("E35").Value = '  -1.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6985'
$ws.Range("E36").Value = '  -3.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.079'
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.180'
$ws.Range("E38").Value = '  -2.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01624'
$ws.Range("E39").Value = '  -7.97%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8742'
$ws.Range("E40").Value = '  -1.79%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4409'
$ws.Range("E41").Value = '  -6.95%  '

$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '104.00'
$ws.Range("E43").Value = '  -4.99%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.743'
$ws.Range("E44").Value = '  -2.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.982'
$ws.Range("E45").Value = '  -8.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.764'
$ws.Range("E46").Value = '  -2.34%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1184'
$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05598'
$ws.Range("E48").Value = '  -4.39%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.66'
$ws.Range("E49").Value = '  -1.36%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.33'
$ws.Range("E50").Value = '  -4.11%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8599'
$ws.Range("E51").Value = '  -1.20%  '
